# Apply updated "dSF" (column F) values for specific rows on Sheet1.
# These values represent a repulled / recalculated data set (see commit message:
# "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = 9
    5  = -3
    6  = 1
    7  = 5
    12 = -7
    13 = 5
    16 = 2
    17 = -3
    18 = 0
    19 = 4
    20 = 1
    21 = 10
    23 = 2
    24 = -8
    25 = 1
    26 = -2
    29 = -3
    30 = 4
    31 = 3
    32 = 10
    34 = -4
    35 = -3
    36 = -3
    37 = 11
    39 = 3
    40 = -4
    41 = -4
    42 = 1
    44 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
